$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "InfoRepeated" instructions for the recognition task (F8/F9):
# removed the blank line before "A döntését így jelölje:" and changed the
# response-key hint from "Új – K" to "Új – J".
$newInstructions = "A döntésre 4 másodperce lesz.`nMinden képet nézzen meg figyelmesen, és minden képre adjon választ, akkor is, ha a döntés nehéz.`nA döntését így jelölje:`nRégi - F`nÚj – J"

$ws.Range("F8").Value = $newInstructions
$ws.Range("F9").Value = $newInstructions

# Move the active selection from F10 to F8, matching the saved view state.
$ws.Range("F8").Select()
